$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (reordered) data for rows 2-10: ID, Spitalname, Lat Geburt, Lon Geburt, Geburtshilfe
$data = @(
    @("SPITAL_Chur",        "Kantonsspital Graubünden, Chur",               "46.861556296686246", "9.542201140210125",  "Ja"),
    @("SPITAL_Samedan",     "Spital Oberengadin, Samedan",                   "46.5398337613234",   "9.878987964240668",  "Ja"),
    @("SPITAL_Ilanz",       "Regionalspital Surselva, Ilanz",                "46.77706472897187",  "9.205079267687351",  "Ja"),
    @("SPITAL_Davos",       "Spital Davos",                                  "46.78780225872632",  "9.814616832894421",  "Ja"),
    @("SPITAL_Schiers",     "Flury Stiftung, Spital Schiers",                "46.97108280571786",  "9.686234727259531",  "Ja"),
    @("SPITAL_Scuol",       "Gesundheitszentrum Unterengadin, Scuol",        "46.7999988668582",   "10.303826533438096", "Ja"),
    @("SPITAL_Bellinzona",  "Ente Ospedaliero Cantonale (EOC), Bellinzona",  "46.184275341809496", "9.026225182498303",  "Ja"),
    @("SPITAL_Thusis",      "Spital Thusis",                                 "46.697267020184185", "9.436968542283159",  "Nein"),
    @("SPITAL_Poschiavo",   "Centro sanitario Valposchiavo",                 "46.32304714423666",  "10.062210768116891", "Nein ")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Drop the "Neonatologie" column (F) entirely - content only, column def/width stays.
$ws.Range("F1:F22").ClearContents()

# Move the active selection, matching the author's final cursor position.
$ws.Range("F6").Select()
